# Trade #46 closed at 2026-02-17 15:30:39 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Summary sheet updates ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.6
$summary.Range("B4").Value = 0.6
$summary.Range("B5").Value = 0.26
$summary.Range("B6").Value = 46
$summary.Range("B8").Value = 23
$summary.Range("B9").Value = 30.43

# --- Strategy Status sheet updates (MarketMaking row) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.6
$status.Range("D4").Value = 46
$status.Range("E4").Value = 0.6
$status.Range("F4").Value = 0.6
$status.Range("G4").Value = 30.43

# --- Append new trade row (Trade #46) to a trade-log sheet ---
# The "Date" column (B) already holds an identical "2026-02-17" string
# elsewhere in the sheet (row 2); copy/paste-special (values only) reuses
# that text verbatim instead of letting plain assignment re-parse the
# date-shaped string into a serial date number.
function Add-TradeRow($ws) {
    $ws.Cells.Item(47, 1).Value = 46
    $ws.Range("B2").Copy()
    $ws.Range("B47").PasteSpecial(-4163)
    $ws.Cells.Item(47, 3).Value = "15:30:33"
    $ws.Cells.Item(47, 4).Value = "MarketMaking"
    $ws.Cells.Item(47, 5).Value = "UP"
    $ws.Cells.Item(47, 6).Value = 0.67
    $ws.Cells.Item(47, 7).Value = 0.5600000000000001
    $ws.Cells.Item(47, 8).Value = "CLOSED"
    $ws.Cells.Item(47, 9).Value = -16.4179
    $ws.Cells.Item(47, 10).Value = -0.11
    $ws.Cells.Item(47, 11).Value = 100.6
    $ws.Cells.Item(47, 12).Value = 0
    $ws.Cells.Item(47, 13).Value = 0
    $ws.Cells.Item(47, 14).Value = 0.6
    $ws.Cells.Item(47, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(47, 16).Value = "early_exit"
    $ws.Cells.Item(47, 17).Value = 0.15
}

# --- All Trades sheet: add row 47 ---
$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades

# --- MarketMaking sheet: add row 47 (duplicate trade log) ---
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $marketMaking
